# Prepare some basic descriptive plots:
# add a "number_of_abstracts" column (G) with counts for years 2004-2020,
# and tidy up a couple of column widths / the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in G1
$ws.Range("G1").Value = "number_of_abstracts"

# New data values in column G for rows 58-74 (years 2004-2020)
$ws.Range("G58").Value = 2002
$ws.Range("G59").Value = 1878
$ws.Range("G60").Value = 2608
$ws.Range("G61").Value = 2526
$ws.Range("G62").Value = 3441
$ws.Range("G63").Value = 2378
$ws.Range("G64").Value = 3305
$ws.Range("G65").Value = 3200
$ws.Range("G66").Value = 3440
$ws.Range("G67").Value = 3255
$ws.Range("G68").Value = 4047
$ws.Range("G69").Value = 4797
$ws.Range("G70").Value = 3426
$ws.Range("G71").Value = 4722
$ws.Range("G72").Value = 4173
$ws.Range("G73").Value = 4464
$ws.Range("G74").Value = 3150

# Resize columns B, C (narrower) and G (new column, wider)
$ws.Columns.Item(2).ColumnWidth = 8.833333333333334
$ws.Columns.Item(3).ColumnWidth = 8.166666666666666
$ws.Columns.Item(7).ColumnWidth = 18.5

# Update the view: scroll down and move the selection
$ws.Range("G78").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 57
$win.ScrollColumn = 1
